# Remove the explicit "VA.MHV.bloodSugarA" example rows (rows 7-8),
# which were a duplicate/explicit variant of the Blood Sugar observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7:K8").EntireRow.Delete()
